# Estado de Cuenta NIT-9003435742
# - Updates the "VALOR MORA" total, worker/period counts
# - Inserts a new worker row (YURI SANTANA GUERRERO) before the signature block,
#   reusing the formatting that previously belonged to the last data row
# - Column E ("Periodo Mora") is center-aligned across the whole table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary updates -------------------------------------------------
$ws.Range("E11").Value = 302505      # VALOR MORA total
$ws.Range("C13").Value = 6           # Cant. Trabajadores
$ws.Range("F13").Value = 9           # Cant. Periodos

# --- Insert a new data row right after the current last row (24) -----------
$ws.Rows.Item(25).Insert()

# Give the new row 25 the special "closing border" look that row 24 used to have
$ws.Range("B24:J24").Copy()
$ws.Range("B25:J25").PasteSpecial(-4122)

# Row 24 is no longer the last row, restore its look to the regular row style
$ws.Range("B23:J23").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Fill the new worker row 25 ---------------------------------------------
$ws.Range("B25").Value = "CC"

$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "1047480370"

$ws.Range("D25").Value = "YURI SANTANA GUERRERO"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2509"

$ws.Range("F25").Value = 56940
$ws.Range("G25").Value = 1423500

# --- Center-align the "Periodo Mora" column across the whole table ---------
$ws.Range("E16:E25").HorizontalAlignment = -4108
